$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns("D").Insert()

# Copy formatting (number format/font/etc.) from the column that used to be D (now E) into the new D column
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Match column width of new column D to column E (both should be width 16)
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth

# Populate the new column D with the FY2018 figures
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 111269000
$ws.Range("D9").Value = 97994000
$ws.Range("D10").Value = 13275000
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 8000
$ws.Range("D15").Value = 1356000
$ws.Range("D17").Value = 106363000
$ws.Range("D18").Value = 4906000
$ws.Range("D20").Value = 3043000
$ws.Range("D21").Value = 9305000
$ws.Range("D22").Value = 504000
$ws.Range("D23").Value = 7445000
$ws.Range("D24").Value = 1536000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 5909000
$ws.Range("D27").Value = 5625000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = -36000
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -3043000
$ws.Range("D33").Value = 5589000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 5589000
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 3019000
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 6173000
$ws.Range("D44").Value = 3543000
$ws.Range("D45").Value = 474000
$ws.Range("D46").Value = 13209000
$ws.Range("D47").Value = 14421000
$ws.Range("D48").Value = 22018000
$ws.Range("D49").Value = 4139000
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 515000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 54302000
$ws.Range("D57").Value = 6113000
$ws.Range("D58").Value = 67000
$ws.Range("D59").Value = 2755000
$ws.Range("D60").Value = 8935000
$ws.Range("D61").Value = 11093000
$ws.Range("D62").Value = 7121000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 29649000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 20489000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 24653000
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 5589000
$ws.Range("D83").Value = 1356000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 7573000
$ws.Range("D91").Value = -2639000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -2471000
$ws.Range("D96").Value = -1436000
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -5167000
$ws.Range("D101").Value = -35000
$ws.Range("D102").Value = -100000
